$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# LOQ4231.xlsx content update (matches the "Introdução à Economia" syllabus
# edit): a new "Objetivos" body paragraph is inserted, several rows shift
# down by one to make room for a "Programa resumido" body + "Programa" body
# + "Método" body, and a brand new "Bibliografia" row (row 22) is appended.
# ---------------------------------------------------------------------------

# Row 10 (Objetivos:) — fill in the body text that was missing before.
$ws.Range("B10").Value = "Apresentar ao aluno de Engenharia conceitos básicos da Ciência Econômica"
$ws.Range("C10").Value = "Apresentar ao aluno de Engenharia conceitos básicos da Ciência Econômica"

# Row 13 used to be "Programa resumido:" / "Semestral" — it becomes the
# (label-less) "5840671 - Francisco José Moreira Chaves" row that used to
# live under "Docentes responsáveis:".
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "5840671 - Francisco José Moreira Chaves"
$ws.Range("C13").Value = "5840671 - Francisco José Moreira Chaves"
$ws.Rows.Item(13).RowHeight = 15

# Row 14 becomes "Programa resumido:" with its real body text (it used to
# be the empty "Short syllabus:" row).
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "História do Pensamento Econômico. Conceitos da Micro e Macroeconomia. Análise da Economia Brasileira"
$ws.Range("C14").Value = "História do Pensamento Econômico. Conceitos da Micro e Macroeconomia. Análise da Economia Brasileira"
$ws.Range("B14").Font.Bold = $false
$ws.Range("B14").WrapText = $true
$ws.Range("C14").Font.Bold = $false
$ws.Range("C14").WrapText = $true
$ws.Rows.Item(14).RowHeight = 60

# Row 15 becomes "Short syllabus:" (label only, used to be "Programa:" with
# a stray "01/01/2012" value that gets cleared).
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Clear()
$ws.Range("C15").Clear()
$ws.Rows.Item(15).RowHeight = 60

# Row 16 becomes "Programa:" with the full program body text (used to be
# the empty "Syllabus:" row).
$ws.Range("A16").Value = "Programa:"
$programaBody = "1.Introdução: história do pensamento econômico.`n2.Microeconomia: oferta, demanda e mercado; elasticidade e estruturas de mercado (concorrência perfeita, monopólio e oligopólio).`n3. Macroeconomia: teoria geral do emprego; juros e a moeda, Sistema Financeiro, Banco Central; Políticas Econômicas: inflação, crescimento, endividamento, balanço de pagamentos e comércio exterior.`n4.Economia brasileira"
$ws.Range("B16").Value = $programaBody
$ws.Range("C16").Value = $programaBody
$ws.Range("B16").Font.Bold = $false
$ws.Range("B16").WrapText = $true
$ws.Range("C16").Font.Bold = $false
$ws.Range("C16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 120

# Row 17 becomes "Syllabus:" (label only, used to be "Avaliação:").
$ws.Range("A17").Value = "Syllabus:"
$ws.Rows.Item(17).RowHeight = 120

# Row 18 becomes "Avaliação:" (label only, used to be "Método:" with a
# stray "5840671 - ..." value that gets cleared).
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows.Item(18).RowHeight = 15

# Row 19 becomes "Método:" with the "Aulas Expositivas..." body (used to be
# "Critério:").
$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras."
$ws.Range("C19").Value = "Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras."
$ws.Rows.Item(19).RowHeight = 60

# Row 20 becomes "Critério:" with the "MF = ..." body (used to be
# "Norma de recuperação:").
$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "MF = (0,30*P1 + 0,60*P2 + 0,20*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários."
$ws.Range("C20").Value = "MF = (0,30*P1 + 0,60*P2 + 0,20*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários."
$ws.Rows.Item(20).RowHeight = 60

# Row 21 becomes "Norma de recuperação:" with the "NF = ..." body (used to
# be "Bibliografia:").
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação."
$ws.Rows.Item(21).RowHeight = 60

# Row 22 is brand new: "Bibliografia:" with the reading list.
$bibliografiaBody = "BEGG, D.; DORNBUSCH, R.; FISCHER, S. Introdução A Economia. Rio de Janeiro: Campus, 2003. `nHUNT, E.K.; SHERMAN, H.J. História do Pensamento Econômico. Petrópolis: Vozes, 2000.`nBACHA , Edmar. Introdução à Macroeconomia: Uma perspectiva brasileira. Rio de Janeiro: Campus,1987.`nROSSETTI, José Pascoal .Introdução à Economia.  9.ed. São Paulo: Atlas, 1982.`nSAMUELSON, P. Introdução à Economia. New York: Mc Graw-Hill Book Company."

$a22 = $ws.Range("A22")
$a22.Value = "Bibliografia:"
$a22.Font.Bold = $true
$a22.WrapText = $false
$a22.VerticalAlignment = -4160

$b22 = $ws.Range("B22")
$b22.Value = $bibliografiaBody
$b22.Font.Bold = $false
$b22.WrapText = $true
$b22.VerticalAlignment = -4160

$c22 = $ws.Range("C22")
$c22.Value = $bibliografiaBody
$c22.Font.Bold = $false
$c22.WrapText = $true
$c22.VerticalAlignment = -4160
$c22.Font.Color = 255

$ws.Rows.Item(22).RowHeight = 120
